# Cost.xlsx: "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to describe a generic "Property" is renamed to
# "DataNode", and the sheet's last active selection moves from A9 to D39
# (the cell that was selected at the moment the author saved the file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the only worksheet: Property1 -> DataNode
$ws.Name = "DataNode"

# Move/record the active selection to D39 (frozen pane stays on bottomLeft)
$ws.Range("D39").Select()
